$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 182, shifting existing
# rows 182-246 down to 183-247 (new dimension becomes A1:R247).
$ws.Rows.Item(182).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = 44900
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = 100112009
$ws.Range("G182").Value = "Acelga"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 40
$ws.Range("K182").Value = 12000
$ws.Range("L182").Value = 12000
$ws.Range("M182").Value = 12000
$ws.Range("N182").Value = "$/docena de atados (12 kilos)"
$ws.Range("O182").Value = "Región de La Araucanía"
$ws.Range("P182").Value = 1000
$ws.Range("Q182").Value = 12
$ws.Range("R182").Value = "Hortaliza"
